$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the species-observation data between row 2 and row 3
# (columns A,B,D,E,F,G,H,Q,R,Z,AB), and moves the "Publik kommentar"
# comment (column AC) from row 2 to row 3. Columns that already held
# identical values in both rows (C, P, S, T, U, V, W, Y, AD, AE, AG,
# AT, AW, AX, AY, etc.) are left untouched.

# --- Row 2 new values (previously row 3's values) ---
$ws.Range("A2").Value2 = 111790625
$ws.Range("B2").Value2 = 96348
$ws.Range("D2").Value2 = "VU"
$ws.Range("E2").Value2 = 220787
$ws.Range("F2").Value2 = "Knärot"
$ws.Range("G2").Value2 = "Goodyera repens"
$ws.Range("H2").Value2 = "(L.) R. Br."
$ws.Range("Q2").Value2 = 489824.6884970492
$ws.Range("R2").Value2 = 6949020.70113107
$ws.Range("Z2").Value2 = "18:29"
$ws.Range("AB2").Value2 = "18:29"
$ws.Range("AC2").Value2 = ""

# --- Row 3 new values (previously row 2's values) ---
$ws.Range("A3").Value2 = 111790785
$ws.Range("B3").Value2 = 77515
$ws.Range("D3").Value2 = "NT"
$ws.Range("E3").Value2 = 6425
$ws.Range("F3").Value2 = "Garnlav"
$ws.Range("G3").Value2 = "Alectoria sarmentosa"
$ws.Range("H3").Value2 = "(Ach.) Ach."
$ws.Range("Q3").Value2 = 489818.2822038208
$ws.Range("R3").Value2 = 6949032.207674611
$ws.Range("Z3").Value2 = "18:34"
$ws.Range("AB3").Value2 = "18:34"
$ws.Range("AC3").Value2 = "Många träd med mycket lav i området"
